$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at U (21) - this shifts existing U:Y to V:Z
$ws.Range("U1").EntireColumn.Insert()

# Header cell for the new "property timezone" column
$ws.Range("U1").Value = "property timezone"
$ws.Range("U1").Style = "Normal"

# Fill data cells U2:U15 with "US/Mountain" and apply style matching siblings (style index 7)
$dataRange = $ws.Range("U2:U15")
$dataRange.Value = "US/Mountain"

# Match styling: copy number format / style attributes from a neighboring data cell that uses style 7 (e.g. G2)
$ws.Range("G2").Copy()
$dataRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Restore sheet view / selection to match target state
$ws.Application.ActiveWindow.ScrollColumn = 16
$ws.Range("U2:U15").Select()
